# Append two more days (2025-12-14 and 2025-12-15) of GSC export data to the
# "Chart" sheet, mirroring the pattern already present for every other row
# (Date in col A, Non-HTTPS count in col B, HTTPS count in col C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper: write $text into $range as literal text (not as a date/number),
# matching how the existing "yyyy-MM-dd" date labels are stored as plain
# shared-string text rather than real date serials.
function Set-TextValue($range, [string]$text) {
    $range.Formula = '=TEXT("' + $text + '","@")'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
    $excel.CutCopyMode = $false
}

Set-TextValue $ws.Range("A70") "2025-12-14"
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 30

Set-TextValue $ws.Range("A71") "2025-12-15"
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = 31
